$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Gabriel"
$ws.Range("B3").Value = 50

$ws.Range("A4").Value = "Clebinho"
$ws.Range("B4").Value = 50

$ws.Range("A5").Value = "Diego"
$ws.Range("B5").Value = 50
